$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 - new issue entry
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "大包客户协议"
$ws.Range("C34").Value = "医院设置里的医院不能删除"
$ws.Range("D34").Value = "未解决"
$ws.Range("E34").Value = 42170
$ws.Range("E34").NumberFormat = "m/d/yy"
$ws.Range("G34").Value = "LiHong"

# Update the view position / selection to match where the user was working
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("E37").Select()
